# Adds NPC interaction-event columns (X:AH) to the results sheet and
# refreshes a handful of previously-rounded numeric values (H, I, U, V, W)
# in rows 2-6 with their more precise counterparts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells (row 1, columns X..AH) -----------------------
# Copy the existing header formatting (bold / border / centered, style
# index 1 in the original file, carried by column W's header cell) onto
# the new header range before writing the header text, so the new
# columns look like the rest of the header row.
$headerRange = $ws.Range("X1:AH1")
$ws.Range("W1").Copy()
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headers = @(
    "NPC_N",
    "NPC_T",
    "NPC_N_type_0",
    "NPC_N_type_1",
    "NPC_N_type_2",
    "NPC_N_actor_1",
    "NPC_N_actor_2",
    "NPC_N_actor_3",
    "NPC_N_actor_4",
    "NPC_N_actor_5",
    "NPC_N_actor_6"
)

$col = 24 # column X
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# --- 2. Refresh a few existing numeric values (rows 2-6) ---------------
# GP_T (H), GP_T_SC (I), OI_T (U), OI_T_SS (V), OI_T_GR (W) were recomputed
# with more precision now that NPC interaction events are concatenated in.
$ws.Range("H2").Value = 57.914
$ws.Range("I2").Value = 49.515
$ws.Range("U2").Value = 8.734
$ws.Range("V2").Value = 0.29
$ws.Range("W2").Value = 19.885

$ws.Range("H3").Value = 47.284
$ws.Range("I3").Value = 45.364
$ws.Range("U3").Value = 3.029
$ws.Range("V3").Value = 0.426

$ws.Range("I4").Value = 35.142
$ws.Range("U4").Value = 1.525
$ws.Range("V4").Value = 0.194
$ws.Range("W4").Value = 5.356

$ws.Range("I5").Value = 68.038
$ws.Range("U5").Value = 6.797
$ws.Range("V5").Value = 0.408
$ws.Range("W5").Value = 21.225

$ws.Range("H6").Value = 43.862
$ws.Range("U6").Value = 10.351
$ws.Range("V6").Value = 0.287
$ws.Range("W6").Value = 47.836

# --- 3. New data cells (rows 2-6, columns X..AH) ------------------------
# Row 2
$ws.Range("X2").Value = 14
$ws.Range("Y2").Value = 76.73999999999999
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 7
$ws.Range("AB2").Value = 7
$ws.Range("AC2").Value = 2
$ws.Range("AD2").Value = 6
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 2
$ws.Range("AH2").Value = 4

# Row 3
$ws.Range("X3").Value = 2
$ws.Range("Y3").Value = 24.092
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 1
$ws.Range("AB3").Value = 1
$ws.Range("AC3").Value = 2
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0

# Row 4
$ws.Range("X4").Value = 6
$ws.Range("Y4").Value = 28.444
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 3
$ws.Range("AB4").Value = 3
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0

# Row 5
$ws.Range("X5").Value = 6
$ws.Range("Y5").Value = 31.68
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = 3
$ws.Range("AB5").Value = 3
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = 2
$ws.Range("AH5").Value = 4

# Row 6 - this scenario had no NPCInteractionEvents at all, so every new
# column is a literal "NULL" (same convention used by I6 already).
$ws.Range("X6").Value = "NULL"
$ws.Range("Y6").Value = "NULL"
$ws.Range("Z6").Value = "NULL"
$ws.Range("AA6").Value = "NULL"
$ws.Range("AB6").Value = "NULL"
$ws.Range("AC6").Value = "NULL"
$ws.Range("AD6").Value = "NULL"
$ws.Range("AE6").Value = "NULL"
$ws.Range("AF6").Value = "NULL"
$ws.Range("AG6").Value = "NULL"
$ws.Range("AH6").Value = "NULL"
